# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp banner
# - Refresh case statistics (Casos totales / Nuevos casos / Casos activos /
#   Recuperados / Casos criticos / Muertes hoy / Muertes) for the countries
#   whose numbers changed in this data pull
# - A handful of neighbouring countries swap places because their refreshed
#   "Casos totales" crossed over their neighbour's count, so the country
#   name in column A is rewritten for those rows as well

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Stats($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 02:00"

# --- Rows whose country stayed in place, only the numbers refreshed ---
Set-Stats 4   6967200 42330 4217939 2545437 0 657 203824   # Estados Unidos
Set-Stats 6   4528347 30913 3820095 571687  0 708 136565   # Brasil
Set-Stats 36  105601  722   80190   23164   0 18  2247     # Panama
Set-Stats 62  48306   2044  24228   23579   0 4   499      # Chequia
Set-Stats 91  14070   48    13365   375     0 1   330      # Zambia
Set-Stats 93  12858   89    10371   2220    0 0   267      # Noruega
Set-Stats 97  10286   55    9681    542     0 0   63       # Guinea
Set-Stats 125 4709    18    4383    229     0 1   97       # Surinam
Set-Stats 154 1904    14    1612    246     0 0   46       # Uruguay

# --- Rows 106/107: Montenegro now ranks above Luxemburgo --------------
$ws.Range("A106").Value = "Montenegro"
Set-Stats 106 7898 187 5129 2635 0 1 134
$ws.Range("A107").Value = "Luxemburgo"
Set-Stats 107 7804 86  6703 977  0 0 124

# --- Rows 130/131: Trinidad yTobago now ranks above Siria --------------
$ws.Range("A130").Value = "Trinidad yTobago"
Set-Stats 130 3853 202 1695 2097 0 1 61
$ws.Range("A131").Value = "Siria"
Set-Stats 131 3765 34  932  2663 0 2 170

# --- Rows 204/205: Santa Lucia now ranks above Timor Oriental ----------
$ws.Range("A204").Value = "Santa Lucia"
Set-Stats 204 27 0 26 1 0 0 0
$ws.Range("A205").Value = "Timor Oriental"
Set-Stats 205 27 0 26 1 0 0 0

# --- Rows 214/215: Montserrat now ranks above Islas Malvinas ------------
$ws.Range("A214").Value = "Montserrat"
Set-Stats 214 13 0 12 0 0 0 1
$ws.Range("A215").Value = "Islas Malvinas"
Set-Stats 215 13 0 13 0 0 0 0
